# Apply the commit "Progess with Dr. Miller's help":
#  - Rename sheets: 1mMGlu -> Glu1mM, 3mMGlu -> Glu3mM
#  - Relabel the time-course headers (J:AB) on every sheet from the old
#    "TimePosition / m1r1..m3r6" scheme to "Position / T0..T59.5"
#  - Update the view state (active sheet / selections) to match the
#    author's last editing position

$wb = $excel.ActiveWorkbook

# --- 1. Rename the glutamate-dose sheets -----------------------------------
$wb.Worksheets.Item("1mMGlu").Name = "Glu1mM"
$wb.Worksheets.Item("3mMGlu").Name = "Glu3mM"

# --- 2. New header labels for columns J..AB (same on every sheet) ----------
$newHeaders = @(
    "Position",
    "T0","T3.5","T7","T10.5","T14","T17.5","T21","T24.5","T28","T31.5",
    "T35","T38.5","T42","T45.5","T49","T52.5","T56","T59.5"
)

foreach ($wsName in @("TimeControl","Glu1mM","Glu3mM")) {
    $ws = $wb.Worksheets.Item($wsName)
    for ($i = 0; $i -lt $newHeaders.Length; $i++) {
        $col = 10 + $i   # column J = 10
        $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
    }
}

# --- 3. Restore the editing/view state --------------------------------------
# TimeControl becomes the active sheet, selection parked on J8
$wsTime = $wb.Worksheets.Item("TimeControl")
[void]$wsTime.Activate()
[void]$wsTime.Range("J8").Select()

# Glu1mM: header row selected (as if the row heading was clicked)
$wsGlu1 = $wb.Worksheets.Item("Glu1mM")
[void]$wsGlu1.Rows(1).Select()

# Glu3mM: single-cell selection left on F24, no longer the active tab
$wsGlu3 = $wb.Worksheets.Item("Glu3mM")
[void]$wsGlu3.Range("F24").Select()

# Make sure TimeControl ends up active/selected last
[void]$wsTime.Activate()
[void]$wsTime.Range("J8").Select()
